$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.454739666666667
$ws.Range("H2").Value = 28.364219
$ws.Range("I2").Value = 0.2644234075555581
$ws.Range("J2").Value = 0.264423407555558
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 59.62013230809134
$ws.Range("R2").Value = 536.581190772822
$ws.Range("S2").Value = 0.00358597939519466
$ws.Range("T2").Value = 0.003585979395194659
$ws.Range("G3").Value = 9.454739666666667
$ws.Range("H3").Value = 28.364219
$ws.Range("I3").Value = 0.2644234075555581
$ws.Range("J3").Value = 0.264423407555558
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 1725.789537380366
$ws.Range("R3").Value = 15532.10583642329
$ws.Range("S3").Value = 0.1038012745343845
$ws.Range("T3").Value = 0.1038012745343845
$ws.Range("G4").Value = 9.454739666666667
$ws.Range("H4").Value = 28.364219
$ws.Range("I4").Value = 0.2644234075555581
$ws.Range("J4").Value = 0.264423407555558
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 1204.499730287356
$ws.Range("R4").Value = 10840.4975725862
$ws.Range("S4").Value = 0.07244719270342492
$ws.Range("T4").Value = 0.07244719270342491
$ws.Range("G5").Value = 9.454739666666667
$ws.Range("H5").Value = 28.364219
$ws.Range("I5").Value = 0.2644234075555581
$ws.Range("J5").Value = 0.264423407555558
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 183.6488758916529
$ws.Range("R5").Value = 1652.839883024876
$ws.Range("S5").Value = 0.01104595141612512
$ws.Range("T5").Value = 0.01104595141612512
$ws.Range("G6").Value = 9.454739666666667
$ws.Range("H6").Value = 28.364219
$ws.Range("I6").Value = 0.2644234075555581
$ws.Range("J6").Value = 0.264423407555558
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 1222.718670102815
$ws.Range("R6").Value = 11004.46803092533
$ws.Range("S6").Value = 0.07354300950642891
$ws.Range("T6").Value = 0.07354300950642888
$ws.Range("I7").Value = 0.4799306150281186
$ws.Range("J7").Value = 0.4799306150281185
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 108.2110204659853
$ws.Range("R7").Value = 973.899184193868
$ws.Range("S7").Value = 0.006508581492553111
$ws.Range("T7").Value = 0.006508581492553109
$ws.Range("I8").Value = 0.4799306150281186
$ws.Range("J8").Value = 0.4799306150281185
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.1884001495500076
$ws.Range("T8").Value = 0.1884001495500075
$ws.Range("I9").Value = 0.4799306150281186
$ws.Range("J9").Value = 0.4799306150281185
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 2186.176714467133
$ws.Range("R9").Value = 19675.59043020419
$ws.Range("S9").Value = 0.1314922384241262
$ws.Range("T9").Value = 0.1314922384241262
$ws.Range("I10").Value = 0.4799306150281186
$ws.Range("J10").Value = 0.4799306150281185
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 333.3241893019048
$ws.Range("R10").Value = 2999.917703717144
$ws.Range("S10").Value = 0.02004849081145659
$ws.Range("T10").Value = 0.02004849081145658
$ws.Range("I11").Value = 0.4799306150281186
$ws.Range("J11").Value = 0.4799306150281185
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 2219.244237012224
$ws.Range("R11").Value = 19973.19813311002
$ws.Range("S11").Value = 0.1334811547499752
$ws.Range("T11").Value = 0.1334811547499751
$ws.Range("G12").Value = 2.798455333333333
$ws.Range("H12").Value = 8.395365999999999
$ws.Range("I12").Value = 0.07826520044130514
$ws.Range("J12").Value = 0.07826520044130512
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 17.64662836987867
$ws.Range("R12").Value = 158.819655328908
$ws.Range("S12").Value = 0.001061393916438094
$ws.Range("T12").Value = 0.001061393916438094
$ws.Range("G13").Value = 2.798455333333333
$ws.Range("H13").Value = 8.395365999999999
$ws.Range("I13").Value = 0.07826520044130514
$ws.Range("J13").Value = 0.07826520044130512
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 510.8067599280224
$ws.Range("R13").Value = 4597.260839352201
$ws.Range("S13").Value = 0.03072355671004506
$ws.Range("T13").Value = 0.03072355671004505
$ws.Range("G14").Value = 2.798455333333333
$ws.Range("H14").Value = 8.395365999999999
$ws.Range("I14").Value = 0.07826520044130514
$ws.Range("J14").Value = 0.07826520044130512
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 356.5131154382793
$ws.Range("R14").Value = 3208.618038944514
$ws.Range("S14").Value = 0.02144323799001064
$ws.Range("T14").Value = 0.02144323799001064
$ws.Range("G15").Value = 2.798455333333333
$ws.Range("H15").Value = 8.395365999999999
$ws.Range("I15").Value = 0.07826520044130514
$ws.Range("J15").Value = 0.07826520044130512
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 54.35720012594044
$ws.Range("R15").Value = 489.214801133464
$ws.Range("S15").Value = 0.003269429169073497
$ws.Range("T15").Value = 0.003269429169073496
$ws.Range("G16").Value = 2.798455333333333
$ws.Range("H16").Value = 8.395365999999999
$ws.Range("I16").Value = 0.07826520044130514
$ws.Range("J16").Value = 0.07826520044130512
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 361.9056371884022
$ws.Range("R16").Value = 3257.15073469562
$ws.Range("S16").Value = 0.02176758265573785
$ws.Range("T16").Value = 0.02176758265573784
$ws.Range("G17").Value = 1.753969666666667
$ws.Range("H17").Value = 5.261909
$ws.Range("I17").Value = 0.04905377116243741
$ws.Range("J17").Value = 0.0490537711624374
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 11.06026260667133
$ws.Range("R17").Value = 99.54236346004201
$ws.Range("S17").Value = 0.0006652429687342824
$ws.Range("T17").Value = 0.0006652429687342822
$ws.Range("G18").Value = 1.753969666666667
$ws.Range("H18").Value = 5.261909
$ws.Range("I18").Value = 0.04905377116243741
$ws.Range("J18").Value = 0.0490537711624374
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 320.1550340183026
$ws.Range("R18").Value = 2881.395306164723
$ws.Range("S18").Value = 0.01925640401676312
$ws.Range("T18").Value = 0.01925640401676311
$ws.Range("G19").Value = 1.753969666666667
$ws.Range("H19").Value = 5.261909
$ws.Range("I19").Value = 0.04905377116243741
$ws.Range("J19").Value = 0.0490537711624374
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 223.4494089647457
$ws.Range("R19").Value = 2011.044680682711
$ws.Range("S19").Value = 0.01343983895029459
$ws.Range("T19").Value = 0.01343983895029459
$ws.Range("G20").Value = 1.753969666666667
$ws.Range("H20").Value = 5.261909
$ws.Range("I20").Value = 0.04905377116243741
$ws.Range("J20").Value = 0.0490537711624374
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 34.06910914395956
$ws.Range("R20").Value = 306.621982295636
$ws.Range("S20").Value = 0.00204915887760109
$ws.Range("T20").Value = 0.00204915887760109
$ws.Range("G21").Value = 1.753969666666667
$ws.Range("H21").Value = 5.261909
$ws.Range("I21").Value = 0.04905377116243741
$ws.Range("J21").Value = 0.0490537711624374
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 226.8292447848478
$ws.Range("R21").Value = 2041.46320306363
$ws.Range("S21").Value = 0.01364312634904433
$ws.Range("T21").Value = 0.01364312634904433
$ws.Range("G22").Value = 4.588468333333334
$ws.Range("H22").Value = 13.765405
$ws.Range("I22").Value = 0.1283270058125809
$ws.Range("J22").Value = 0.1283270058125809
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 28.93417468587667
$ws.Range("R22").Value = 260.40757217289
$ws.Range("S22").Value = 0.001740307346255843
$ws.Range("T22").Value = 0.001740307346255842
$ws.Range("G23").Value = 4.588468333333334
$ws.Range("H23").Value = 13.765405
$ws.Range("I23").Value = 0.1283270058125809
$ws.Range("J23").Value = 0.1283270058125809
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 837.5408442165597
$ws.Range("R23").Value = 7537.867597949036
$ws.Range("S23").Value = 0.05037567166866078
$ws.Range("T23").Value = 0.05037567166866076
$ws.Range("G24").Value = 4.588468333333334
$ws.Range("H24").Value = 13.765405
$ws.Range("I24").Value = 0.1283270058125809
$ws.Range("J24").Value = 0.1283270058125809
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 584.5543150613885
$ws.Range("R24").Value = 5260.988835552495
$ws.Range("S24").Value = 0.03515925993505018
$ws.Range("T24").Value = 0.03515925993505017
$ws.Range("G25").Value = 4.588468333333334
$ws.Range("H25").Value = 13.765405
$ws.Range("I25").Value = 0.1283270058125809
$ws.Range("J25").Value = 0.1283270058125809
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 89.12641502462445
$ws.Range("R25").Value = 802.1377352216201
$ws.Range("S25").Value = 0.005360697393194076
$ws.Range("T25").Value = 0.005360697393194074
$ws.Range("G26").Value = 4.588468333333334
$ws.Range("H26").Value = 13.765405
$ws.Range("I26").Value = 0.1283270058125809
$ws.Range("J26").Value = 0.1283270058125809
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 593.3961268253723
$ws.Range("R26").Value = 5340.565141428351
$ws.Range("S26").Value = 0.02970638548048269
$ws.Range("T26").Value = 0.02970638548048269
